$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# The row for account 004290978 / LARISSA / 5397.89 was removed from the
# export (duplicate/erroneous entry). Locate it by its account number and
# delete the whole row, shifting the following rows up by one.
$target = $ws.Cells.Find("004290978")
if ($target -ne $null) {
    $target.EntireRow.Delete()
} else {
    # Fallback: known position in the original export (row 6 - the 5th
    # data row, right after the header row).
    $ws.Rows(6).Delete()
}
